# Daily attendance processing - 2026-01-29 23:12:49
# Reorders the "Recorded By" (column G) contributor lists on sheet1:
#   "dnasr281@gmail.com, System"          -> "System, dnasr281@gmail.com"
#   "system, backup@backdoor.com, System" -> "backup@backdoor.com, System, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Text

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "system, backup@backdoor.com, System") {
        $cell.Value = "backup@backdoor.com, System, system"
    }
}
